$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.23672555474949
$ws.Range("D2").Value = 8.967251477545098
$ws.Range("E2").Value = 14.01533542610956
$ws.Range("F2").Value = 34.48826144649518
$ws.Range("G2").Value = 3.666989002550985
$ws.Range("I2").Value = 27.25246170100857
$ws.Range("J2").Value = 10.54649845162567
$ws.Range("K2").Value = 10.44298049628918
$ws.Range("L2").Value = 10.34638984231915
$ws.Range("N2").Value = 20.29987755906077
$ws.Range("O2").Value = 26.1098727306632

$ws.Range("B3").Value = 14.11317823058214
$ws.Range("D3").Value = 8.970603754623804
$ws.Range("E3").Value = 14.0443047989054
$ws.Range("F3").Value = 34.5560757195231
$ws.Range("G3").Value = 3.668778695742112
$ws.Range("I3").Value = 27.35228124614035
$ws.Range("J3").Value = 10.56319261295698
$ws.Range("K3").Value = 10.21792051646982
$ws.Range("L3").Value = 10.31479671890644
$ws.Range("N3").Value = 20.35700642307954
$ws.Range("O3").Value = 26.17576803935464

$ws.Range("B4").Value = 14.03904805365975
$ws.Range("D4").Value = 8.973667341000587
$ws.Range("E4").Value = 14.06335448924256
$ws.Range("F4").Value = 34.60460789323179
$ws.Range("G4").Value = 3.669937060898431
$ws.Range("I4").Value = 27.41749553706517
$ws.Range("J4").Value = 10.57400776149459
$ws.Range("K4").Value = 10.0786856767158
$ws.Range("L4").Value = 10.29678191591949
$ws.Range("N4").Value = 20.39377507221801
$ws.Range("O4").Value = 26.2210937188578

$ws.Range("B5").Value = 14.00930252456749
$ws.Range("D5").Value = 8.975169338371064
$ws.Range("E5").Value = 14.07143547055514
$ws.Range("F5").Value = 34.62611694461557
$ws.Range("G5").Value = 3.670424108173057
$ws.Range("I5").Value = 27.44505877512615
$ws.Range("J5").Value = 10.57855744958474
$ws.Range("K5").Value = 10.02176177651072
$ws.Range("L5").Value = 10.28979367334398
$ws.Range("N5").Value = 20.40918508657609
$ws.Range("O5").Value = 26.24078650700488

$ws.Range("B6").Value = 14.00439209075949
$ws.Range("D6").Value = 8.975434080115392
$ws.Range("E6").Value = 14.07279654137198
$ws.Range("F6").Value = 34.62979304504942
$ws.Range("G6").Value = 3.670505889557687
$ws.Range("I6").Value = 27.44969532521538
$ws.Range("J6").Value = 10.57932153548252
$ws.Range("K6").Value = 10.01230080953013
$ws.Range("L6").Value = 10.28865473877376
$ws.Range("N6").Value = 20.41176970212846
$ws.Range("O6").Value = 26.24413025695732

$ws.Range("B7").Value = 14.03864498283001
$ws.Range("D7").Value = 8.973686569802835
$ws.Range("E7").Value = 14.06346218343291
$ws.Range("F7").Value = 34.60489096239485
$ws.Range("G7").Value = 3.669943568576259
$ws.Range("I7").Value = 27.4178632632912
$ws.Range("J7").Value = 10.57406854297837
$ws.Range("K7").Value = 10.07791862250215
$ws.Range("L7").Value = 10.29668623457901
$ws.Range("N7").Value = 20.39398116861858
$ws.Range("O7").Value = 26.22135435585368

$ws.Range("B8").Value = 14.19378504385623
$ws.Range("D8").Value = 8.968199176564445
$ws.Range("E8").Value = 14.02506244584103
$ws.Range("F8").Value = 34.51021202095595
$ws.Range("G8").Value = 3.667593769960383
$ws.Range("I8").Value = 27.2860655063395
$ws.Range("J8").Value = 10.55213759392025
$ws.Range("K8").Value = 10.36564353842566
$ws.Range("L8").Value = 10.33521250156796
$ws.Range("N8").Value = 20.31922528875413
$ws.Range("O8").Value = 26.13158280444996

$ws.Range("B9").Value = 14.51044269033805
$ws.Range("D9").Value = 8.96537953337285
$ws.Range("E9").Value = 13.95974829855418
$ws.Range("F9").Value = 34.37931358189019
$ws.Range("G9").Value = 3.66345573320217
$ws.Range("I9").Value = 27.05870607105241
$ws.Range("J9").Value = 10.51359535138636
$ws.Range("K9").Value = 10.91825023469468
$ws.Range("L9").Value = 10.42150604521359
$ws.Range("N9").Value = 20.1859933105383
$ws.Range("O9").Value = 25.99420315220478

$ws.Range("B10").Value = 14.74890966985767
$ws.Range("D10").Value = 8.968099923833172
$ws.Range("E10").Value = 13.91781118878745
$ws.Range("F10").Value = 34.31660325726774
$ws.Range("G10").Value = 3.66069907093279
$ws.Range("I10").Value = 26.9105589653655
$ws.Range("J10").Value = 10.48797515320533
$ws.Range("K10").Value = 11.31296828862819
$ws.Range("L10").Value = 10.49112504189588
$ws.Range("N10").Value = 20.09617593983311
$ws.Range("O10").Value = 25.91690929084186

$ws.Range("B11").Value = 14.85827746997332
$ws.Range("D11").Value = 8.970366559375334
$ws.Range("E11").Value = 13.90003807527352
$ws.Range("F11").Value = 34.29534854598297
$ws.Range("G11").Value = 3.659505945010992
$ws.Range("I11").Value = 26.84725281004311
$ws.Range("J11").Value = 10.47690018765414
$ws.Range("K11").Value = 11.48927402468421
$ws.Range("L11").Value = 10.52407090050863
$ws.Range("N11").Value = 20.05705081547175
$ws.Range("O11").Value = 25.88689144994906

$ws.Range("B12").Value = 14.8997878291698
$ws.Range("D12").Value = 8.97137178466407
$ws.Range("E12").Value = 13.89349476954694
$ws.Range("F12").Value = 34.28834564323945
$ws.Range("G12").Value = 3.659062848773769
$ws.Range("I12").Value = 26.82386723182178
$ws.Range("J12").Value = 10.47278936460895
$ws.Range("K12").Value = 11.55550651644599
$ws.Range("L12").Value = 10.53672325099678
$ws.Range("N12").Value = 20.04248315059651
$ws.Range("O12").Value = 25.87626484464921

$ws.Range("B13").Value = 14.89084411535997
$ws.Range("D13").Value = 8.971148775901709
$ws.Range("E13").Value = 13.89489568024863
$ws.Range("F13").Value = 34.28980733600621
$ws.Range("G13").Value = 3.659157890550766
$ws.Range("I13").Value = 26.82887763238018
$ws.Range("J13").Value = 10.47367101679241
$ws.Range("K13").Value = 11.54126671972917
$ws.Range("L13").Value = 10.53399060623404
$ws.Range("N13").Value = 20.04560953730825
$ws.Range("O13").Value = 25.87852052830756

$ws.Range("B14").Value = 14.86169083116541
$ws.Range("D14").Value = 8.970446322196128
$ws.Range("E14").Value = 13.89949600932625
$ws.Range("F14").Value = 34.29475145512401
$ws.Range("G14").Value = 3.659469316800104
$ws.Range("I14").Value = 26.84531710148736
$ws.Range("J14").Value = 10.47656032577917
$ws.Range("K14").Value = 11.49473399131244
$ws.Range("L14").Value = 10.52510832117455
$ws.Range("N14").Value = 20.0558473572841
$ws.Range("O14").Value = 25.88600234765442

$ws.Range("B15").Value = 14.84384503048369
$ws.Range("D15").Value = 8.970035146500571
$ws.Range("E15").Value = 13.90233817808537
$ws.Range("F15").Value = 34.29791605728307
$ws.Range("G15").Value = 3.659661208042067
$ws.Range("I15").Value = 26.85546319080506
$ws.Range("J15").Value = 10.47834091315618
$ws.Range("K15").Value = 11.46616045320332
$ws.Range("L15").Value = 10.51969043505952
$ws.Range("N15").Value = 20.06215060526262
$ws.Range("O15").Value = 25.89068163579973

$ws.Range("B16").Value = 14.7417773104553
$ws.Range("D16").Value = 8.967972410754568
$ws.Range("E16").Value = 13.91899890050142
$ws.Range("F16").Value = 34.31813865558397
$ws.Range("G16").Value = 3.660778266272945
$ws.Range("I16").Value = 26.91477835661884
$ws.Range("J16").Value = 10.48871056537327
$ws.Range("K16").Value = 11.30137535696817
$ws.Range("L16").Value = 10.48899703297563
$ws.Range("N16").Value = 20.0987676445063
$ws.Range("O16").Value = 25.91897460613979

$ws.Range("B17").Value = 14.67936631265335
$ws.Range("D17").Value = 8.966969815511058
$ws.Range("E17").Value = 13.92955336803716
$ws.Range("F17").Value = 34.33240738056265
$ws.Range("G17").Value = 3.6614791117889
$ws.Range("I17").Value = 26.95221257399345
$ws.Range("J17").Value = 10.49522026144001
$ws.Range("K17").Value = 11.19940567398531
$ws.Range("L17").Value = 10.47048949411973
$ws.Range("N17").Value = 20.12167417452119
$ws.Range("O17").Value = 25.9376494768736

$ws.Range("B18").Value = 14.64355468203413
$ws.Range("D18").Value = 8.966490083071337
$ws.Range("E18").Value = 13.93574681817652
$ws.Range("F18").Value = 34.34129892487377
$ws.Range("G18").Value = 3.661887953658605
$ws.Range("I18").Value = 26.97412843812372
$ws.Range("J18").Value = 10.49901905968079
$ws.Range("K18").Value = 11.14045265902537
$ws.Range("L18").Value = 10.45996496818322
$ws.Range("N18").Value = 20.13501260934049
$ws.Range("O18").Value = 25.94887485998343

$ws.Range("B19").Value = 14.63144516168926
$ws.Range("D19").Value = 8.966344334558091
$ws.Range("E19").Value = 13.9378649222852
$ws.Range("F19").Value = 34.34442701384764
$ws.Range("G19").Value = 3.662027366600292
$ws.Range("I19").Value = 26.98161486249029
$ws.Range("J19").Value = 10.50031465456587
$ws.Range("K19").Value = 11.12044216845253
$ws.Range("L19").Value = 10.45642245336216
$ws.Range("N19").Value = 20.13955683748983
$ws.Range("O19").Value = 25.95275869585242

$ws.Range("B20").Value = 14.68600145459436
$ws.Range("D20").Value = 8.967066518898353
$ws.Range("E20").Value = 13.92841712179209
$ws.Range("F20").Value = 34.33081760048232
$ws.Range("G20").Value = 3.661403912436714
$ws.Range("I20").Value = 26.94818782769038
$ws.Range("J20").Value = 10.49452164551351
$ws.Range("K20").Value = 11.21029232311772
$ws.Range("L20").Value = 10.4724472296839
$ws.Range("N20").Value = 20.11921885200349
$ws.Range("O20").Value = 25.93561139653804

$ws.Range("B21").Value = 14.87025153373206
$ws.Range("D21").Value = 8.97064867199672
$ws.Range("E21").Value = 13.89813971091369
$ws.Range("F21").Value = 34.29327086694436
$ws.Range("G21").Value = 3.659377607179734
$ws.Range("I21").Value = 26.8404725019052
$ws.Range("J21").Value = 10.47570941513943
$ws.Range("K21").Value = 11.50841666277519
$ws.Range("L21").Value = 10.5277125311672
$ws.Range("N21").Value = 20.05283353456537
$ws.Range("O21").Value = 25.88378465293018

$ws.Range("B22").Value = 14.99120924701577
$ws.Range("D22").Value = 8.973845516279223
$ws.Range("E22").Value = 13.8794413123689
$ws.Range("F22").Value = 34.27482729966754
$ws.Range("G22").Value = 3.658104076106921
$ws.Range("I22").Value = 26.77349606954862
$ws.Range("J22").Value = 10.46389831570059
$ws.Range("K22").Value = 11.70013896579253
$ws.Range("L22").Value = 10.56485679406882
$ws.Range("N22").Value = 20.01089293283754
$ws.Range("O22").Value = 25.85422936022916

$ws.Range("B23").Value = 14.92661316589489
$ws.Range("D23").Value = 8.972061376878832
$ws.Range("E23").Value = 13.88932148440833
$ws.Range("F23").Value = 34.28411335724812
$ws.Range("G23").Value = 3.658779151257582
$ws.Range("I23").Value = 26.8089297420823
$ws.Range("J23").Value = 10.47015796964605
$ws.Range("K23").Value = 11.59811820005712
$ws.Range("L23").Value = 10.54494076544865
$ws.Range("N23").Value = 20.03314545837976
$ws.Range("O23").Value = 25.86960837533104

$ws.Range("B24").Value = 14.6830014886306
$ws.Range("D24").Value = 8.967022498105063
$ws.Range("E24").Value = 13.92893042762753
$ws.Range("F24").Value = 34.33153419617031
$ws.Range("G24").Value = 3.661437891631652
$ws.Range("I24").Value = 26.95000618681379
$ws.Range("J24").Value = 10.49483731442253
$ws.Range("K24").Value = 11.20537148986668
$ws.Range("L24").Value = 10.47156177649872
$ws.Range("N24").Value = 20.12032837628616
$ws.Range("O24").Value = 25.93653128966769

$ws.Range("B25").Value = 14.42362782289627
$ws.Range("D25").Value = 8.965297030830088
$ws.Range("E25").Value = 13.97635236142272
$ws.Range("F25").Value = 34.40885281608128
$ws.Range("G25").Value = 3.664525177429988
$ws.Range("I25").Value = 27.11689044381743
$ws.Range("J25").Value = 10.52354668871092
$ws.Range("K25").Value = 10.77044749283085
$ws.Range("L25").Value = 10.39704415633026
$ws.Range("N25").Value = 20.22061357800364
$ws.Range("O25").Value = 26.02722119989309
